$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: bold font, thin border all around, centered horizontally, top-aligned
# vertically.
$b1 = $ws.Range("B1")
$b1.Font.Bold = $true
$b1.Borders.LineStyle = 1
$b1.VerticalAlignment = -4160
$b1.HorizontalAlignment = -4108
$b1.Value = 0

# A2 gets the exact same formatting as B1 - copy it over rather than
# re-applying property by property so both cells share one style record.
$b1.Copy()
$a2 = $ws.Range("A2")
$a2.PasteSpecial(-4122)
$a2.Value = 0

# B2 = "disconnected_elements" (plain, no special style)
$ws.Range("B2").Value = "disconnected_elements"
